$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cell = $ws.Cells.Item(554, 5)
$url = "https://www.instagram.com/p/DKzwzpmsycB/?img_index=1&igsh=MW1ia2pmMWYyNDloaQ=="
$cell.Value = $url
$chars = $cell.Characters(1, $url.Length)
$chars.Font.Name = "Calibri"
$chars.Font.Size = 11
$chars.Font.Underline = $true
$chars.Font.ColorIndex = 11
Write-Host ("style before: " + $cell.Style)
$cell.NumberFormat = "@"
Write-Host ("style after numfmt: " + $cell.Style)
